$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# fix calculations and data scaling: the diameter measurement columns
# (B:I) were recorded in the wrong units - scale rows 2-12 up by 10.
$range = $ws.Range("B2:I12")
foreach ($cell in $range.Cells) {
    $v = $cell.Value()
    $cell.Value = $v * 10
}

# Update the selected/active cell to match the saved view state.
$ws.Range("L13").Select()
